$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row 5 data - shared string insertion order matters, so set
# the description (G) first, then id (A), then name (E), matching the
# order new unique strings appear in the target file.
$ws.Range("G5").Value = "a dataset with a lot of variables"
$ws.Range("A5").Value = "dataset_4"
$ws.Range("E5").Value = "dataset with lot of variables"
$ws.Range("F5").Value = "open_data"
$ws.Range("H5").Value = 200
$ws.Range("K5").NumberFormat = "@"
$ws.Range("L5").NumberFormat = "@"

# Expand the table to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P5"))

# Adjust column widths (name/description/data_path columns got wider to
# fit the new long strings)
$ws.Columns.Item(5).ColumnWidth = 15.166666666666666
$ws.Columns.Item(7).ColumnWidth = 51.5
$ws.Columns.Item(16).ColumnWidth = 22.0

# Update selection / scroll position
$ws.Range("E6").Select()
